# Add a new trailing date column (P) to the attendance sheet:
#   - P1 holds the new date header "2025-04-06" (kept as literal text, like
#     the other date headers in row 1, rather than an auto-converted date
#     serial number).
#   - P2:P8 hold FALSE (not yet marked present), matching the pattern of the
#     other attendance columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the new header cell so "2025-04-06" isn't
# auto-parsed into a date serial, then drop back to the default/Normal
# style so no stray per-cell number format lingers on the sheet.
$ws.Range("P1").NumberFormat = "@"
$ws.Range("P1").Value = "2025-04-06"
$ws.Range("P1").Style = "Normal"

# New attendance column starts fully unchecked.
$ws.Range("P2").Value = $false
$ws.Range("P3").Value = $false
$ws.Range("P4").Value = $false
$ws.Range("P5").Value = $false
$ws.Range("P6").Value = $false
$ws.Range("P7").Value = $false
$ws.Range("P8").Value = $false
